$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.678.66'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '2.623.54'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.87'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.20'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  -0.93%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -2.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.80'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.65'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '3.097.39'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000184'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '63.603.65'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '2.631.76'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.10'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  +0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.50'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '343.15'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.01'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.88'
$ws.Range('E24').Value = '  +7.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000112'
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('E27').Value = '  -0.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '572.04'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.25'
$ws.Range('E29').Value = '  +4.57%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.04'
$ws.Range('E32').Value = '  -1.86%  '
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.65'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').Value = '  +1.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.403'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.68'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.90'
$ws.Range('E39').Value = '  +0.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '153.33'
$ws.Range('E40').Value = '  -0.77%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.54'
$ws.Range('E42').Value = '  +5.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '163.30'
$ws.Range('E43').Value = '  +3.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.08'
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.91'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0585'
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.629'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0999'
$ws.Range('E48').Value = '  -2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0248'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('D50').Value = '0.0₆0235'
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.791'
$ws.Range('E51').Value = '  +1.43%  '
